$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "G2" = 3.770298333333333
    "H2" = 11.310895
    "I2" = 0.06100259562224731
    "J2" = 0.06125631726190612
    "K2" = 2
    "L2" = 0.6666666666666666
    "M2" = 0.9763746666666667
    "N2" = 2.929124
    "O2" = 0.1257320070262716
    "P2" = 0.1269157938307497
    "Q2" = 3.681223778442222
    "R2" = 33.13101400598001
    "S2" = 0.007669978781397206
    "T2" = 0.007774394132443071
    "G3" = 3.770298333333333
    "H3" = 11.310895
    "I3" = 0.06100259562224731
    "J3" = 0.06125631726190612
    "O3" = 0.2098967719105039
    "P3" = 0.2118729833364321
    "Q3" = 6.145427930803334
    "R3" = 55.30885137723
    "S3" = 0.01280424789927155
    "T3" = 0.01297855868648303
    "G4" = 3.770298333333333
    "H4" = 11.310895
    "I4" = 0.06100259562224731
    "J4" = 0.06125631726190612
    "M4" = 2.486016
    "N4" = 7.458048
    "O4" = 0.3201350791356976
    "P4" = 0.3231492017230527
    "Q4" = 9.373021981439999
    "R4" = 84.35719783296
    "S4" = 0.0195290707770111
    "T4" = 0.01979493002367902
    "G5" = 3.770298333333333
    "H5" = 11.310895
    "I5" = 0.06100259562224731
    "J5" = 0.06125631726190612
    "K5" = 1
    "L5" = 0.5
    "M5" = 0.217295
    "N5" = 0.43459
    "O5" = 0.02798202104121269
    "P5" = 0.0188303174740658
    "Q5" = 0.8192669763416666
    "R5" = 4.91560185805
    "S5" = 0.001706975914270313
    "T5" = 0.001153475901333789
    "G6" = 3.770298333333333
    "H6" = 11.310895
    "I6" = 0.06100259562224731
    "J6" = 0.06125631726190612
    "K6" = 3
    "L6" = 1
    "M6" = 2.455878333333333
    "N6" = 7.367635
    "O6" = 0.3162541208863143
    "P6" = 0.3192317036356997
    "Q6" = 9.259393987036111
    "R6" = 83.33454588332501
    "S6" = 0.01929232225029715
    "T6" = 0.01955495851796721
    "I7" = 0.06469423882843597
    "J7" = 0.06496331472897099
    "K7" = 2
    "L7" = 0.6666666666666666
    "M7" = 0.9763746666666667
    "N7" = 2.929124
    "O7" = 0.1257320070262716
    "P7" = 0.1269157938307497
    "Q7" = 3.903997327887556
    "R7" = 35.135975950988
    "S7" = 0.008134136490936207
    "T7" = 0.008244870658704186
    "I8" = 0.06469423882843597
    "J8" = 0.06496331472897099
    "O8" = 0.2098967719105039
    "P8" = 0.2118729833364321
    "S8" = 0.01357911189129589
    "T8" = 0.01376397129905066
    "I9" = 0.06469423882843597
    "J9" = 0.06496331472897099
    "M9" = 2.486016
    "N9" = 7.458048
    "O9" = 0.3201350791356976
    "P9" = 0.3231492017230527
    "Q9" = 9.940241336064
    "R9" = 89.46217202457601
    "S9" = 0.02071089526696507
    "T9" = 0.02099284329595041
    "I10" = 0.06469423882843597
    "J10" = 0.06496331472897099
    "K10" = 1
    "L10" = 0.5
    "M10" = 0.217295
    "N10" = 0.43459
    "O10" = 0.02798202104121269
    "P10" = 0.0188303174740658
    "Q10" = 0.8688458727216667
    "R10" = 5.21307523633
    "S10" = 0.001810275552142534
    "T10" = 0.001223279840514178
    "I11" = 0.06469423882843597
    "J11" = 0.06496331472897099
    "K11" = 3
    "L11" = 1
    "M11" = 2.455878333333333
    "N11" = 7.367635
    "O11" = 0.3162541208863143
    "P11" = 0.3192317036356997
    "Q11" = 9.819737011082779
    "R11" = 88.37763309974501
    "S11" = 0.02045981962709628
    "T11" = 0.02073834963475156
    "G12" = 31.40746233333333
    "H12" = 94.222387
    "I12" = 0.5081658147055464
    "J12" = 0.5102793749960634
    "K12" = 2
    "L12" = 0.6666666666666666
    "M12" = 0.9763746666666667
    "N12" = 2.929124
    "O12" = 0.1257320070262716
    "P12" = 0.1269157938307497
    "Q12" = 30.66545056655422
    "R12" = 275.989055098988
    "S12" = 0.06389270778506881
    "T12" = 0.06476251195308419
    "G13" = 31.40746233333333
    "H13" = 94.222387
    "I13" = 0.5081658147055464
    "J13" = 0.5102793749960634
    "O13" = 0.2098967719105039
    "P13" = 0.2118729833364321
    "Q13" = 51.19284448991533
    "R13" = 460.735600409238
    "S13" = 0.1066623641019655
    "T13" = 0.1081144135154659
    "G14" = 31.40746233333333
    "H14" = 94.222387
    "I14" = 0.5081658147055464
    "J14" = 0.5102793749960634
    "M14" = 2.486016
    "N14" = 7.458048
    "O14" = 0.3201350791356976
    "P14" = 0.3231492017230527
    "Q14" = 78.07945388006399
    "R14" = 702.7150849205759
    "S14" = 0.1626817033048163
    "T14" = 0.1648963726857161
    "G15" = 31.40746233333333
    "H15" = 94.222387
    "I15" = 0.5081658147055464
    "J15" = 0.5102793749960634
    "K15" = 1
    "L15" = 0.5
    "M15" = 0.217295
    "N15" = 0.43459
    "O15" = 0.02798202104121269
    "P15" = 0.0188303174740658
    "Q15" = 6.824684527721666
    "R15" = 40.94810716633
    "S15" = 0.01421950651951559
    "T15" = 0.009608722631643744
    "G16" = 31.40746233333333
    "H16" = 94.222387
    "I16" = 0.5081658147055464
    "J16" = 0.5102793749960634
    "K16" = 3
    "L16" = 1
    "M16" = 2.455878333333333
    "N16" = 7.367635
    "O16" = 0.3162541208863143
    "P16" = 0.3192317036356997
    "Q16" = 77.1329062494161
    "R16" = 694.196156244745
    "S16" = 0.1607095329941803
    "T16" = 0.1628973542101534
    "G17" = 0.7679895
    "H17" = 1.535979
    "I17" = 0.01242590075603175
    "J17" = 0.008318388326620067
    "K17" = 2
    "L17" = 0.6666666666666666
    "M17" = 0.9763746666666667
    "N17" = 2.929124
    "O17" = 0.1257320070262716
    "P17" = 0.1269157938307497
    "Q17" = 0.749845492066
    "R17" = 4.499072952396
    "S17" = 0.001562333441165137
    "T17" = 0.001055734857865427
    "G18" = 0.7679895
    "H18" = 1.535979
    "I18" = 0.01242590075603175
    "J18" = 0.008318388326620067
    "O18" = 0.2098967719105039
    "P18" = 0.2118729833364321
    "Q18" = 1.251790629441
    "R18" = 7.510743776646
    "S18" = 0.002608156456771354
    "T18" = 0.001762441751311945
    "G19" = 0.7679895
    "H19" = 1.535979
    "I19" = 0.01242590075603175
    "J19" = 0.008318388326620067
    "M19" = 2.486016
    "N19" = 7.458048
    "O19" = 0.3201350791356976
    "P19" = 0.3231492017230527
    "Q19" = 1.909234184832
    "R19" = 11.455405108992
    "S19" = 0.003977966721864547
    "T19" = 0.002688080547369635
    "G20" = 0.7679895
    "H20" = 1.535979
    "I20" = 0.01242590075603175
    "J20" = 0.008318388326620067
    "K20" = 1
    "L20" = 0.5
    "M20" = 0.217295
    "N20" = 0.43459
    "O20" = 0.02798202104121269
    "P20" = 0.0188303174740658
    "Q20" = 0.1668802784025
    "R20" = 0.66752111361
    "S20" = 0.0003477018164113009
    "T20" = 0.0001566378930628188
    "G21" = 0.7679895
    "H21" = 1.535979
    "I21" = 0.01242590075603175
    "J21" = 0.008318388326620067
    "K21" = 3
    "L21" = 1
    "M21" = 2.455878333333333
    "N21" = 7.367635
    "O21" = 0.3162541208863143
    "P21" = 0.3192317036356997
    "Q21" = 1.8860887732775
    "R21" = 11.316532639665
    "S21" = 0.003929742319819409
    "T21" = 0.002655493277010242
    "G22" = 21.861327
    "H22" = 65.58398100000001
    "I22" = 0.3537114500877385
    "J22" = 0.3551826046864394
    "K22" = 2
    "L22" = 0.6666666666666666
    "M22" = 0.9763746666666667
    "N22" = 2.929124
    "O22" = 0.1257320070262716
    "P22" = 0.1269157938307497
    "Q22" = 21.344845862516
    "R22" = 192.103612762644
    "S22" = 0.04447285052770426
    "T22" = 0.04507828222865281
    "G23" = 21.861327
    "H23" = 65.58398100000001
    "I23" = 0.3537114500877385
    "J23" = 0.3551826046864394
    "O23" = 0.2098967719105039
    "P23" = 0.2118729833364321
    "Q23" = 35.633044834266
    "R23" = 320.697403508394
    "S23" = 0.07424289156119965
    "T23" = 0.07525359808412051
    "G24" = 21.861327
    "H24" = 65.58398100000001
    "I24" = 0.3537114500877385
    "J24" = 0.3551826046864394
    "M24" = 2.486016
    "N24" = 7.458048
    "O24" = 0.3201350791356976
    "P24" = 0.3231492017230527
    "Q24" = 54.34760870323201
    "R24" = 489.1284783290881
    "S24" = 0.1132354430650405
    "T24" = 0.1147769751703375
    "G25" = 21.861327
    "H25" = 65.58398100000001
    "I25" = 0.3537114500877385
    "J25" = 0.3551826046864394
    "K25" = 1
    "L25" = 0.5
    "M25" = 0.217295
    "N25" = 0.43459
    "O25" = 0.02798202104121269
    "P25" = 0.0188303174740658
    "Q25" = 4.750357050465
    "R25" = 28.50214230279
    "S25" = 0.00989756123887295
    "T25" = 0.006688201207511263
    "G26" = 21.861327
    "H26" = 65.58398100000001
    "I26" = 0.3537114500877385
    "J26" = 0.3551826046864394
    "K26" = 3
    "L26" = 1
    "M26" = 2.455878333333333
    "N26" = 7.367635
    "O26" = 0.3162541208863143
    "P26" = 0.3192317036356997
    "Q26" = 53.688759317215
    "R26" = 483.1988338549351
    "S26" = 0.1118627036949212
    "T26" = 0.1133855479958173
}

foreach ($cellRef in $changes.Keys) {
    $ws.Range($cellRef).Value = $changes[$cellRef]
}

Write-Host "Updated $($changes.Count) cells"
